# Updated symbol list on Fri Jan  6 03:59:07 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "257.91"
Set-TextValue "E2" "0.02%"
Set-TextValue "D3" "27.01"
Set-TextValue "E3" "-0.33%"
Set-TextValue "D4" "4.642"
Set-TextValue "E4" "-5.32%"
Set-TextValue "D5" "0.05911"
Set-TextValue "D6" "6.647"
Set-TextValue "E6" "-0.60%"
Set-TextValue "D7" "0.8555"
Set-TextValue "E7" "-1.51%"
Set-TextValue "D8" "0.9500"
Set-TextValue "E8" "-1.18%"
Set-TextValue "D9" "0.1399"
Set-TextValue "E9" "-0.93%"
Set-TextValue "D10" "0.05201"
Set-TextValue "E10" "46.39%"
Set-TextValue "D11" "0.07096"
Set-TextValue "E11" "-1.36%"
Set-TextValue "D12" "0.03107"
Set-TextValue "E12" "-1.13%"
Set-TextValue "D13" "0.09139"
Set-TextValue "E13" "-1.15%"
Set-TextValue "D14" "0.001539"
Set-TextValue "E14" "-0.45%"
Set-TextValue "D15" "0.01052"
Set-TextValue "E15" "1,630.49%"
Set-TextValue "D16" "0.006157"
Set-TextValue "E16" "2.80%"
Set-TextValue "D17" "3.498"
Set-TextValue "E17" "0.40%"
Set-TextValue "D18" "3.189"
Set-TextValue "E18" "-0.80%"
Set-TextValue "E19" "-1.02%"
Set-TextValue "D20" "0.3056"
Set-TextValue "E20" "-2.84%"
Set-TextValue "E21" "-2.21%"
Set-TextValue "D22" "3.822"
Set-TextValue "E22" "8.29%"
Set-TextValue "D23" "0.04276"
Set-TextValue "E23" "0.09%"
Set-TextValue "D24" "0.001219"
Set-TextValue "E24" "-0.23%"
Set-TextValue "D25" "0.004296"
Set-TextValue "E25" "-4.92%"
Set-TextValue "E26" "-0.04%"
Set-TextValue "E27" "29.84%"
Set-TextValue "D40" "0.03829"
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006208"
Set-TextValue "E41" "-5.66%"
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1102"
Set-TextValue "E42" "-0.09%"
Set-TextValue "D43" "0.002199"
Set-TextValue "E43" "-0.04%"
Set-TextValue "E44" "32.27%"
Set-TextValue "D45" "0.00005000"
Set-TextValue "E45" "-8.92%"
Set-TextValue "E46" "-0.04%"
Set-TextValue "E47" "-53.27%"
Set-TextValue "D48" "0.2495"
Set-TextValue "E48" "11,621.88%"
Set-TextValue "E49" "-0.04%"
Set-TextValue "E50" "-0.04%"
